# Add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the "总计" (Total) sheet,
#    and populate it with the per-fund holdings detail for the quarter.
# 2. Update the "总计" (Total) summary sheet with a new leading row for
#    2022-Q1 (existing rows shift down by one).

$wb = $excel.ActiveWorkbook

# "总计" is the last sheet before this edit (position 6).
$totalSheetBeforeInsert = $wb.Worksheets.Item(6)

# Reference cells carrying the two formats used throughout this workbook:
#   - $styledCell: bold, centered, thin-bordered -- used for header row and
#     the running-index column A.
#   - $blankCell:  completely unstyled -- used to reset formatting on data
#     cells after a NumberFormat="@" (Text) trick is used to preserve
#     leading zeros / exact decimal text such as "008297" or "0.0060".
$styledCell = $wb.Worksheets.Item(5).Range("A2")
$blankCell = $wb.Worksheets.Item(5).Range("C2")

# ---------------------------------------------------------------------
# 1. New sheet "2022-Q1" (inserted before "总计")
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($totalSheetBeforeInsert)
$newSheet.Name = "2022-Q1"

# NOTE: inserting a sheet "before" $totalSheetBeforeInsert re-targets that
# positional handle at the freshly-inserted sheet, not the original "总计"
# tab (which has now shifted one slot further down). Re-resolve "总计" by
# name so later edits land on the right worksheet/XML part.
$totalSheet = $wb.Worksheets.Item("总计")

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"
$styledCell.Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

function Set-FundRow($row, $idx, $code, $name, $scale, $stockPos, $ratio, $heldValue, $rank) {
    $newSheet.Range("A$row").Value = $idx
    $styledCell.Copy()
    $newSheet.Range("A$row").PasteSpecial(-4122)

    $newSheet.Range("B$row`:G$row").NumberFormat = "@"
    $newSheet.Range("B$row").Value = $code
    $newSheet.Range("C$row").Value = $name
    $newSheet.Range("D$row").Value = $scale
    $newSheet.Range("E$row").Value = $stockPos
    $newSheet.Range("F$row").Value = $ratio
    $newSheet.Range("G$row").Value = $heldValue
    $blankCell.Copy()
    $newSheet.Range("B$row`:G$row").PasteSpecial(-4122)

    $newSheet.Range("H$row").Value = $rank
}

Set-FundRow 2 0 "008297" "广发价值优势混合"         "28.25" "93.96" "5.24" "1.4803" 9
Set-FundRow 3 1 "270022" "广发内需增长混合A"         "15.92" "79.56" "5.04" "0.8024" 10
Set-FundRow 4 2 "011134" "广发价值优选混合A"         "6.35"  "93.95" "5.56" "0.3531" 6
Set-FundRow 5 3 "011135" "广发价值优选混合C"         "1.48"  "93.95" "5.56" "0.0823" 6
Set-FundRow 6 4 "005043" "国寿安保健康科学混合A"     "0.99"  "85.72" "2.65" "0.0262" 9
Set-FundRow 7 5 "005044" "国寿安保健康科学混合C"     "0.87"  "85.72" "2.65" "0.0231" 9
Set-FundRow 8 6 "011183" "广发内需增长混合C"         "0.12"  "79.56" "5.04" "0.0060" 10
Set-FundRow 9 7 "003366" "浙商汇金中证转型成长指数"   "0.09"  "93.88" "1.19" "0.0011" 10

# ---------------------------------------------------------------------
# 2. Rebuild "总计" with the new 2022-Q1 row on top
# ---------------------------------------------------------------------
$totalSheet.Cells.Clear()

$totalSheet.Range("B1").Value = "日期"
$totalSheet.Range("C1").Value = "持有数量(只)"
$totalSheet.Range("D1").Value = "持有市值(亿元)"
$styledCell.Copy()
$totalSheet.Range("B1:D1").PasteSpecial(-4122)

function Set-TotalRow($row, $idx, $date, $count, $heldValue) {
    $totalSheet.Range("A$row").Value = $idx
    $styledCell.Copy()
    $totalSheet.Range("A$row").PasteSpecial(-4122)
    $totalSheet.Range("B$row").Value = $date
    $totalSheet.Range("C$row").Value = $count
    $totalSheet.Range("D$row").Value = $heldValue
}

Set-TotalRow 2 0 "2022-Q1" 8 2.77
Set-TotalRow 3 1 "2021-Q4" 1 0.11
Set-TotalRow 4 2 "2021-Q3" 9 0.4
Set-TotalRow 5 3 "2021-Q2" 19 2.66
Set-TotalRow 6 4 "2021-Q1" 16 2.17
Set-TotalRow 7 5 "2020-Q4" 9 2.28
